$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Fix the typo "defininição" -> "definição" ("softwares terceiros, a
# defininição de padrões").  "defininição" = "defini" + "ni" + "ção",
# so the fix simply drops the duplicated "ni".
# ------------------------------------------------------------------

# Locate the misspelled phrase through to the end of its run (the
# trailing space before "ou") on a throw-away duplicate range so the
# document's real Find/Selection state is left alone.
$findRng = $d.Content.Duplicate
$findRng.Find.Execute("defininição de padrões ", $true, $false, $false, `
                       $false, $false, $true, 1, $false, "", 0)

$wordStart = $findRng.Start
$runEnd    = $findRng.End

# Split point: right after "defini", i.e. where the stray "ni" starts.
$splitPoint = $wordStart + 6
$dupEnd     = $splitPoint + 2

# A throw-away bookmark placed at the enclosing run's natural end (the
# boundary with the next run, "ou") keeps that boundary from being
# swallowed into the following text once the duplicate letters are
# deleted a moment later.
$barrierRng = $d.Range($runEnd, $runEnd)
$d.Bookmarks.Add("ZZZ_TempBarrier", $barrierRng)

# Word itself re-drops its "_GoBack" (last-edit-location) bookmark right
# at the point being corrected; re-adding a bookmark under that name
# simply relocates it, exactly like Word does after an edit.
$editRng = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $editRng)

# Remove the duplicated "ni".
$dupRng = $d.Range($splitPoint, $dupEnd)
$dupRng.Text = ""

# Drop the temporary barrier bookmark now that the edit is complete.
$d.Bookmarks("ZZZ_TempBarrier").Delete()
